# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 15 (pushing the previous row 15 down
# to row 16, which keeps all of its original values/format), then fill
# the new row 15 with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 15; Excel shifts row 15 (and
# everything below it) down by one, so the old row 15 becomes row 16.
$ws.Rows(15).Insert()

# Populate the newly inserted row 15 with the new weekly record.
$ws.Range("A15").Value = 10
$ws.Range("B15").Value = "Vega Modelo de Temuco"
$ws.Range("C15").Value = "La Araucanía"
$ws.Range("D15").Value = 44753
$ws.Range("E15").Value = 9
$ws.Range("F15").Value = "Fruta"
$ws.Range("G15").Value = 100108
$ws.Range("H15").Value = "Tropicales y subtropicales"
$ws.Range("I15").Value = 100108001
$ws.Range("J15").Value = "Guayaba"
$ws.Range("K15").Value = "Sin especificar"
$ws.Range("L15").Value = "Primera"
$ws.Range("M15").Value = 160
$ws.Range("N15").Value = 2300
$ws.Range("O15").Value = 2300
$ws.Range("P15").Value = 2300
$ws.Range("Q15").Value = "$/kilo"
$ws.Range("R15").Value = "Región de Arica y Parinacota"
$ws.Range("S15").Value = 2300
$ws.Range("T15").Value = 1
